$d = $word.ActiveDocument

# Update the two date stamp cells: 07/13/2015 -> 07/21/2015
$d.Content.Find.Execute("07/13/2015", $false, $false, $false, $false, $false,
                         $true, 1, $false, "07/21/2015", 2)

# Update the three "Estimated Completion Date" notes: July 2 2015 -> July 3 2015
$d.Content.Find.Execute("July 2 2015", $false, $false, $false, $false, $false,
                         $true, 1, $false, "July 3 2015", 2)
